# Actualización automática hashcode jue jun 27 01:42:02 CEST 2019
# Update specific hashcode values (column B) in the hashcode.csv sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hashcode.csv")

$updates = @{
    "B100" = "aed0b10bec2e9f418199ec1ba1e7d195"
    "B104" = "d3250a5c1d6eed791df60eeb4e0dcd1e"
    "B113" = "d2c1e61c41dbe7d42161673c62f26d94"
    "B122" = "6403db4eaca423e88668dcf035f15b05"
    "B164" = "2c5adae7a14dbf122a3e7e333946951c"
    "B230" = "6cc38a03a89a547d65027e64ed10d11e"
    "B233" = "bf8a8bb894e8abfbc617dca6f524925d"
    "B331" = "e40c86b9f34ec9b35c007636ab13d92c"
    "B342" = "091253ce18fb2eaae2c3e52d7191a868"
    "B343" = "9683867abb9ed08c897898b9ce16f688"
    "B381" = "426758b07b194188b97fe09b886f440d"
    "B458" = "62f05aaa5756711c583f9c74bdffd409"
    "B477" = "e1b8840a7130774ea1c4a2335241f85b"
    "B619" = "5f08ecffecd63e81c0870b802b54b76b"
    "B623" = "ff9f888e91bca8d85efafc7661513a32"
    "B628" = "846c9647ded4ae397a5a92e7ec1d0301"
    "B757" = "d9fc2c056853e5b1d706e3bee6e41e35"
    "B760" = "3abdca00a5fbc0b0fc986be33a67f946"
    "B763" = "2664bab3d905a7c9c5fbeecc3e104f3c"
    "B767" = "9cae9b93fd2c3bed47e712b319e14674"
    "B779" = "d7f4356c35eb2b8b0deaac7d4e0be00c"
    "B818" = "1dcbd17e31672161575ab11d3dad0626"
    "B831" = "572bd04638f6b5cead7fe5e2de230d72"
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
